$d = $word.ActiveDocument

$replacements = @(
    @("96×42=4032", "96×60=5760"),
    @("87×56=4872", "41×28=1148"),
    @("61×22=1342", "15×99=1485"),
    @("28×90=2520", "70×52=3640"),
    @("35×11=385",  "21×19=399"),
    @("59×43=2537", "97×16=1552"),
    @("95×53=5035", "60×87=5220"),
    @("28×49=1372", "87×41=3567"),
    @("17×92=1564", "86×73=6278"),
    @("74×61=4514", "87×49=4263"),
    @("84×32=2688", "58×55=3190"),
    @("56×17=952",  "65×55=3575"),
    @("73×73=5329", "55×67=3685"),
    @("48×21=1008", "62×55=3410"),
    @("31×31=961",  "90×73=6570"),
    @("75×30=2250", "77×57=4389"),
    @("35×70=2450", "63×99=6237"),
    @("27×50=1350", "62×48=2976"),
    @("65×78=5070", "58×51=2958"),
    @("46×62=2852", "46×93=4278"),
    @("80×68=5440", "98×27=2646"),
    @("52×52=2704", "69×98=6762"),
    @("68×25=1700", "55×48=2640"),
    @("13×76=988",  "11×61=671"),
    @("65×49=3185", "60×99=5940")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
